$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.823918879032135
$ws.Range("B1").Value = 2.668290615081787
$ws.Range("C1").Value = 6.797377586364746
$ws.Range("D1").Value = 4.806951999664307
$ws.Range("E1").Value = 2.431468963623047
